$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'65.248.21"
$ws.Range('E2').Value = '  -3.87%  '

$ws.Range('D3').Value = "'3.395.21"
$ws.Range('E3').Value = '  -5.32%  '

$ws.Range('E4').Value = '  +0.36%  '

$ws.Range('D5').Value = "'185.53"
$ws.Range('E5').Value = '  -8.21%  '

$ws.Range('D6').Value = "'530.44"
$ws.Range('E6').Value = '  -5.54%  '

$ws.Range('D7').Value = "'0.609"
$ws.Range('E7').Value = '  -1.52%  '

$ws.Range('D8').Value = "'3.386.29"
$ws.Range('E8').Value = '  -5.49%  '

$ws.Range('E9').Value = '  +0.04%  '

$ws.Range('D10').Value = "'0.630"
$ws.Range('E10').Value = '  -5.59%  '

$ws.Range('D11').Value = "'58.54"
$ws.Range('E11').Value = '  -2.47%  '

$ws.Range('D12').Value = "'0.134"
$ws.Range('E12').Value = '  -10.56%  '

$ws.Range('D13').Value = "'0.0000257"
$ws.Range('E13').Value = '  -9.45%  '

$ws.Range('D14').Value = "'9.31"
$ws.Range('E14').Value = '  -6.36%  '

$ws.Range('D15').Value = "'3.961.41"
$ws.Range('E15').Value = '  -4.66%  '

$ws.Range('E16').Value = '  -2.48%  '

$ws.Range('D17').Value = "'3.412.12"
$ws.Range('E17').Value = '  -4.81%  '

$ws.Range('D18').Value = "'65.155.83"
$ws.Range('E18').Value = '  -3.67%  '

$ws.Range('D19').Value = "'17.58"
$ws.Range('E19').Value = '  -6.74%  '

$ws.Range('D20').Value = "'11.27"
$ws.Range('E20').Value = '  -8.36%  '

$ws.Range('D21').Value = "'0.978"
$ws.Range('E21').Value = '  -8.50%  '

$ws.Range('D22').Value = "'375.10"
$ws.Range('E22').Value = '  -6.31%  '

$ws.Range('D23').Value = "'82.24"
$ws.Range('E23').Value = '  -3.08%  '

$ws.Range('D24').Value = "'3.75"
$ws.Range('E24').Value = '  -9.16%  '

$ws.Range('D25').Value = "'10.94"
$ws.Range('E25').Value = '  -14.43%  '

$ws.Range('D26').Value = "'3.70"
$ws.Range('E26').Value = '  -5.17%  '

$ws.Range('D27').Value = "'11.72"
$ws.Range('E27').Value = '  -6.37%  '

$ws.Range('D28').Value = "'2.67"
$ws.Range('E28').Value = '  -7.72%  '

$ws.Range('D29').Value = "'8.57"
$ws.Range('E29').Value = '  -7.70%  '

$ws.Range('D30').Value = "'680.87"
$ws.Range('E30').Value = '  +2.24%  '

$ws.Range('D31').Value = "'29.82"
$ws.Range('E31').Value = '  -5.09%  '

$ws.Range('D32').Value = "'6.81"
$ws.Range('E32').Value = '  -17.15%  '

$ws.Range('D33').Value = "'11.26"
$ws.Range('E33').Value = '  -7.22%  '

$ws.Range('D34').Value = "'61.42"
$ws.Range('E34').Value = '  -3.01%  '

$ws.Range('D35').Value = "'0.106"
$ws.Range('E35').Value = '  -6.22%  '

$ws.Range('D36').Value = "'0.999"
$ws.Range('E36').Value = '  -0.27%  '

$ws.Range('D37').Value = "'36.68"
$ws.Range('E37').Value = '  -12.37%  '

$ws.Range('D38').Value = "'0.387"
$ws.Range('E38').Value = '  -7.70%  '

$ws.Range('D39').Value = "'1.00"
$ws.Range('E39').Value = '  +0.47%  '

$ws.Range('D40').Value = "'0.128"
$ws.Range('E40').Value = '  -5.05%  '

$ws.Range('D41').Value = "'28.97"
$ws.Range('E41').Value = '  +30.28%  '

$ws.Range('D42').Value = "'2.895.32"
$ws.Range('E42').Value = '  -12.09%  '

$ws.Range('D43').Value = "'2.80"
$ws.Range('E43').Value = '  -11.35%  '

$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = "'2.67"
$ws.Range('E44').Value = '  -3.00%  '

$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = "'0.0397"
$ws.Range('E45').Value = '  -4.22%  '

$ws.Range('B46').Value = 'PEPE'
$ws.Range('C46').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D46').Value = "'0.0₃0629"
$ws.Range('E46').Value = '  -16.87%  '

$ws.Range('D47').Value = "'2.38"
$ws.Range('E47').Value = '  -13.36%  '

$ws.Range('D48').Value = "'0.126"
$ws.Range('E48').Value = '  -3.18%  '

$ws.Range('D49').Value = "'136.96"
$ws.Range('E49').Value = '  -1.52%  '

$ws.Range('D50').Value = "'2.90"
$ws.Range('E50').Value = '  -6.94%  '

$ws.Range('D51').Value = "'2.64"
$ws.Range('E51').Value = '  -2.62%  '

# Restore default (unstyled) cell formatting on the Price cells we just
# wrote as text, so only the cell VALUES differ from the original file.
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
